# Inserts a new weekly price record for "Camote" / "1a (guarda)" at row 918
# on the "Terminal La Palmera de La Serena - Zapallo" sheet, shifting all
# subsequent rows down by one (dimension grows from A1:R966 to A1:R967).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 918..966 down to 919..967, creating a blank row 918.
$ws.Rows.Item(918).Insert()

# Populate the new row 918 with the new data point.
$row = 918
$ws.Cells.Item($row, 1).Value  = 8
$ws.Cells.Item($row, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item($row, 3).Value  = 'Coquimbo'
$ws.Cells.Item($row, 4).Value  = 45147
$ws.Cells.Item($row, 5).Value  = 4
$ws.Cells.Item($row, 6).Value  = 100112045
$ws.Cells.Item($row, 7).Value  = 'Zapallo'
$ws.Cells.Item($row, 8).Value  = 'Camote'
$ws.Cells.Item($row, 9).Value  = '1a (guarda)'
$ws.Cells.Item($row, 10).Value = 1800
$ws.Cells.Item($row, 11).Value = 700
$ws.Cells.Item($row, 12).Value = 800
$ws.Cells.Item($row, 13).Value = 750
$ws.Cells.Item($row, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item($row, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item($row, 16).Value = 750
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = 'Hortaliza'

# Match the date formatting used by the rest of column D (yyyy-mm-dd hh:mm:ss).
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row + 1, 4).NumberFormat
